# Auto-generated edit script: updates crypto price/volume table
# matching commit "Updated cryptos list on Mon Apr 10 22:57:58 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''29.689.59'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +3.79%  '
$ws.Range('D3').Value = '''1.909.68'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.54%  '
$ws.Range('D4').Value = '''1.003'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.64%  '
$ws.Range('D5').Value = '''316.07'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.10%  '
$ws.Range('E6').Value = '  -0.60%  '
$ws.Range('D7').Value = '''0.5168'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.13%  '
$ws.Range('E8').Value = '  +0.77%  '
$ws.Range('D9').Value = '''0.08497'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.83%  '
$ws.Range('D10').Value = '''42.68'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.26%  '
$ws.Range('D11').Value = '''1.122'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.60%  '
$ws.Range('D12').Value = '''6.308'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.25%  '
$ws.Range('D13').Value = '''1.912.46'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.69%  '
$ws.Range('D14').Value = '''20.93'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.95%  '
$ws.Range('D15').Value = '''7.351'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.83%  '
$ws.Range('D16').Value = '''1.002'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.66%  '
$ws.Range('D17').Value = '''93.41'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.02%  '
$ws.Range('D18').Value = '''0.00001118'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.77%  '
$ws.Range('D19').Value = '''0.06757'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.43%  '
$ws.Range('D20').Value = '''17.95'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.99%  '
$ws.Range('D21').Value = '''1.001'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.57%  '
$ws.Range('D22').Value = '''6.039'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.91%  '
$ws.Range('D23').Value = '''29.702.25'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +3.74%  '
$ws.Range('D24').Value = '''11.23'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.66%  '
$ws.Range('D25').Value = '''2.211'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.88%  '
$ws.Range('D26').Value = '''2.124.13'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.33%  '
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').Value = '''20.98'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.83%  '
$ws.Range('B28').Value = 'Monero'
$ws.Range('C28').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D28').Value = '''159.20'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.22%  '
$ws.Range('D29').Value = '''2.449'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.76%  '
$ws.Range('D30').Value = '''128.22'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.07%  '
$ws.Range('D31').Value = '''1.077'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.86%  '
$ws.Range('D32').Value = '''0.1054'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.13%  '
$ws.Range('E33').Value = '  +5.96%  '
$ws.Range('D34').Value = '''3.667'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.38%  '
$ws.Range('D35').Value = '''0.02497'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.40%  '
$ws.Range('D36').Value = '''0.06629'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.18%  '
$ws.Range('D37').Value = '''9.090'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.30%  '
$ws.Range('D38').Value = '''0.2203'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.58%  '
$ws.Range('B39').Value = 'ARBITRUM'
$ws.Range('C39').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D39').Value = '''1.238'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +3.13%  '
$ws.Range('B40').Value = 'InternetComputer(DFINITY)'
$ws.Range('C40').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D40').Value = '''5.217'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.37%  '
$ws.Range('D41').Value = '''0.6554'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.88%  '
$ws.Range('D42').Value = '''1.240'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.81%  '
$ws.Range('D43').Value = '''11.33'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.94%  '
$ws.Range('D44').Value = '''0.6115'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.45%  '
$ws.Range('D45').Value = '''13.18'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.06%  '
$ws.Range('D46').Value = '''3.679'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.68%  '
$ws.Range('D47').Value = '''2.065'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.14%  '
$ws.Range('D48').Value = '''1.238'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.39%  '
$ws.Range('D49').Value = '''124.26'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.23%  '
$ws.Range('D50').Value = '''1.160'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.32%  '
$ws.Range('D51').Value = '''78.28'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.17%  '
